# case6.xlsx edit: "add ep factory"
#
# 1. Merge A2:B2 and put a thin bottom border under the merged header cell.
# 2. Add a new label "调入单位主管领导签字：" (with a trailing tab) as wrapped,
#    centered text in D4, matching the existing bordered-cell look used
#    elsewhere on the sheet.
# 3. Grow row 4 to fit the new two-line label.
# 4. Leave the selection parked on D4 (where the edit was made).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A2:B2 merge + bottom border -------------------------------------------------
$ws.Range("A2:B2").Merge()
# The merge copies A2's style onto B2; put B2 back to the workbook default
# formatting (no explicit font) before the border is applied, same as the
# rest of the sheet's "blank" cells.
$ws.Range("B2").Style = "Normal"
$ws.Range("A2:B2").Borders.Item(9).LineStyle = 1

# --- D4 label ----------------------------------------------------------------
# D4 already has the standard bordered/centered data-cell style used across
# row 4-8; just add the text and turn wrapping on so the two-line label fits.
$ws.Range("D4").Value = "调入单位主管领导签字：	"
$ws.Range("D4").WrapText = $true

# --- Row height ---------------------------------------------------------------
$ws.Rows.Item(4).RowHeight = 30

# --- Selection ------------------------------------------------------------------
$ws.Range("D4").Select()
